$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin rows (prices/volume refreshed, coin ranking shifted down for
# rows 7-23, and the "Hora" run counter bumped from 0 to 2) per the Jan 27 2023
# GitHub Actions refresh. Columns D/E/G are forced to Text format cell-by-cell
# before assignment so numeric-looking values (prices, "%", counters) are kept
# as literal text, matching the original inlineStr cells.
# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '302.43'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-1.95%'
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = '2'
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '35.60'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.22%'
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = '2'
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.014'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-2.53%'
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = '2'
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07867'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-3.61%'
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = '2'
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.850'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-4.67%'
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = '2'
# Row 7
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '7.810'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '0.30%'
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = '2'
# Row 8
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9248'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '-0.82%'
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = '2'
# Row 9
$ws.Range("B9").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C9").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1317'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-5.60%'
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = '2'
# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1886'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-2.63%'
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = '2'
# Row 11
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.09110'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-1.63%'
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = '2'
# Row 12
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03513'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '-1.22%'
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = '2'
# Row 13
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09831'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.34%'
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = '2'
# Row 14
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001409'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-0.71%'
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = '2'
# Row 15
$ws.Range("B15").Value = 'CoinExToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.04426'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-1.46%'
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = '2'
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006122'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '3.93%'
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = '2'
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.656'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '3.28%'
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = '2'
# Row 18
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.150'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-0.52%'
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = '2'
# Row 19
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.265'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '9.33%'
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = '2'
# Row 20
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3442'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '0.05%'
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = '2'
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1333'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.25%'
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = '2'
# Row 22
$ws.Range("B22").Value = 'MCDex'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.212'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '6.63%'
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = '2'
# Row 23
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.2514'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '4.69%'
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = '2'
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.001238'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '2.28%'
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = '2'
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.004619'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-5.89%'
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = '2'
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0001305'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '5.31%'
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = '2'
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0004459'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '0.40%'
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = '2'
# Row 28
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = '2'
# Row 29
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = '2'
# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = '2'
# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = '2'
# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = '2'
# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = '2'
# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = '2'
# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = '2'
# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = '2'
# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = '2'
# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = '2'
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01921'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '-4.26%'
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = '2'
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05075'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '2.85%'
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = '2'
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007559'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '-1.07%'
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = '2'
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.01019'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '-8.60%'
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = '2'
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1343'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-2.99%'
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = '2'
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.002168'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '3.32%'
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = '2'
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.009949'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-6.27%'
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = '2'
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006143'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '-4.75%'
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = '2'
# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.45%'
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = '2'
# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '63.57'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-1.69%'
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = '2'
# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '40.02%'
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = '2'
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002108'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.45%'
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = '2'
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002007'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.45%'
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = '2'
